# CompStat 90th Precinct weekly report refresh (new crime data collected).
# Updates the report header (volume/week-of text) and the weekly numeric
# table (rows 15-30) to the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Header text: "Volume 30   Number  44" -> "...  45"
#              "Report Covering the Week  10/30/2023  Through  11/5/2023"
#              -> "...11/6/2023  Through  11/12/2023"
# ---------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 30   Number  45"
$ws.Range("C9").Value = "Report Covering the Week  11/6/2023  Through  11/12/2023"

# ---------------------------------------------------------------------
# Weekly crime-stat table: plain numeric updates (style/type unchanged)
# ---------------------------------------------------------------------
$ws.Range("M15").Value = 37.5

$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 13
$ws.Range("H16").Value = 30.769230769230
$ws.Range("I16").Value = 166
$ws.Range("J16").Value = 154
$ws.Range("K16").Value = 7.792207792207
$ws.Range("L16").Value = 44.347826086956
$ws.Range("M16").Value = -49.848942598187
$ws.Range("N16").Value = -85.231316725978

$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 11
$ws.Range("E17").Value = -81.818181818181
$ws.Range("F17").Value = 19
$ws.Range("G17").Value = 25
$ws.Range("H17").Value = -24
$ws.Range("I17").Value = 263
$ws.Range("J17").Value = 250
$ws.Range("K17").Value = 5.2
$ws.Range("L17").Value = 35.567010309278
$ws.Range("M17").Value = 47.752808988764
$ws.Range("N17").Value = -49.423076923076

$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -40
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -34.375
$ws.Range("I18").Value = 245
$ws.Range("J18").Value = 268
$ws.Range("K18").Value = -8.582089552238
$ws.Range("L18").Value = 24.365482233502
$ws.Range("M18").Value = -40.243902439024
$ws.Range("N18").Value = -79.237288135593

$ws.Range("C19").Value = 13
$ws.Range("D19").Value = 13
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 67
$ws.Range("G19").Value = 66
$ws.Range("H19").Value = 1.515151515151
$ws.Range("I19").Value = 675
$ws.Range("J19").Value = 595
$ws.Range("K19").Value = 13.445378151260
$ws.Range("L19").Value = 53.758542141230
$ws.Range("M19").Value = 54.462242562929
$ws.Range("N19").Value = 42.405063291139

$ws.Range("G20").Value = 14
$ws.Range("H20").Value = -21.428571428571
$ws.Range("I20").Value = 157
$ws.Range("J20").Value = 160
$ws.Range("K20").Value = -1.875
$ws.Range("L20").Value = 16.296296296296
$ws.Range("M20").Value = 12.949640287769
$ws.Range("N20").Value = -80.448318804483

$ws.Range("C21").Value = 29
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -30.952380952381
$ws.Range("F21").Value = 135
$ws.Range("G21").Value = 150
$ws.Range("H21").Value = -10
$ws.Range("I21").Value = 1522
$ws.Range("J21").Value = 1442
$ws.Range("K21").Value = 5.547850208044
$ws.Range("L21").Value = 38.237965485921
$ws.Range("M21").Value = 1.129568106312
$ws.Range("N21").Value = -63.369434416365

$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 20
$ws.Range("K22").Value = -44.444444444444
$ws.Range("M22").Value = -48.717948717948

$ws.Range("D23").Value = 7
$ws.Range("E23").Value = -71.428571428571
$ws.Range("F23").Value = 11
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = -42.105263157894
$ws.Range("I23").Value = 172
$ws.Range("J23").Value = 142
$ws.Range("K23").Value = 21.126760563380
$ws.Range("L23").Value = 17.006802721088
$ws.Range("M23").Value = 42.148760330578

$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 17
$ws.Range("E24").Value = -17.647058823529
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 82
$ws.Range("H24").Value = -26.829268292682
$ws.Range("I24").Value = 903
$ws.Range("J24").Value = 1058
$ws.Range("K24").Value = -14.650283553875
$ws.Range("L24").Value = 4.878048780487
$ws.Range("M24").Value = -16.155988857938

$ws.Range("C25").Value = 8
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = -33.333333333333
$ws.Range("G25").Value = 47
$ws.Range("H25").Value = -12.765957446808
$ws.Range("I25").Value = 433
$ws.Range("J25").Value = 437
$ws.Range("K25").Value = -0.915331807780
$ws.Range("L25").Value = 30.815709969788
$ws.Range("M25").Value = 1.643192488262

$ws.Range("F27").Value = 5
$ws.Range("H27").Value = -37.5
$ws.Range("I27").Value = 49
$ws.Range("K27").Value = 11.363636363636
$ws.Range("L27").Value = -22.222222222222

$ws.Range("N28").Value = -84.810126582278

$ws.Range("N29").Value = -87.671232876712

$ws.Range("G30").Value = 4
$ws.Range("H30").Value = -50
$ws.Range("J30").Value = 29
$ws.Range("K30").Value = -44.827586206896

# ---------------------------------------------------------------------
# Cells that flip between a numeric value and the sheet's literal
# placeholder text ("0" / "***.*") used when a category has no prior
# week data. Assigning a bare numeric-looking string auto-coerces to a
# real number, so force text with a leading apostrophe, then re-pull
# the cell formatting (number format / alignment / quote-prefix state)
# from a same-row donor cell that already carries the right look, via
# a formats-only paste so the style index matches cells of the same
# kind elsewhere in the table.
# ---------------------------------------------------------------------

# Row 22: "2023" (C) and "28-day prior" (D)/"28-day %chg" (E) swap kind.
$ws.Range("C22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("C22").PasteSpecial(-4122)

$ws.Range("D22").Value = "'0"
$ws.Range("A22").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Range("E22").Value = "'***.*"
$ws.Range("A22").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F22").Value = 1
$ws.Range("G22").Copy()
$ws.Range("F22").PasteSpecial(-4122)

# Row 27: same pattern.
$ws.Range("C27").Value = 1
$ws.Range("F27").Copy()
$ws.Range("C27").PasteSpecial(-4122)

$ws.Range("D27").Value = "'0"
$ws.Range("A27").Copy()
$ws.Range("D27").PasteSpecial(-4122)

$ws.Range("E27").Value = "'***.*"
$ws.Range("A27").Copy()
$ws.Range("E27").PasteSpecial(-4122)

# Row 30: same pattern, but C30 becomes text while D30/E30 become numeric.
$ws.Range("C30").Value = "'0"
$ws.Range("A30").Copy()
$ws.Range("C30").PasteSpecial(-4122)

$ws.Range("D30").Value = 1
$ws.Range("F30").Copy()
$ws.Range("D30").PasteSpecial(-4122)

$ws.Range("E30").Value = -100
$ws.Range("H30").Copy()
$ws.Range("E30").PasteSpecial(-4122)

Write-Output "CompStat weekly refresh applied"
